$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "IsAvailable" header (D1) to "Availability" and turn on wrap text,
# matching the new shared-string entry + the new wrapText cell style.
$ws.Range("D1").Value = "Availability"
$ws.Range("D1").WrapText = $true

# The whole D column (rows 2-107) used to hold boolean TRUE; it now holds the
# plain number 0.
$ws.Range("D2:D107").Value = 0

# Give column D an explicit width (it previously had none) and move the
# active selection to D6, matching the saved view state. (22.8333 is the
# COM ColumnWidth value that round-trips to the closest persisted column
# width to the recorded 23.7109375 character-width units.)
$ws.Columns("D").ColumnWidth = 22.8333
$ws.Range("D6").Select() | Out-Null
